# Adding full run for ZEV Jan R2-4 and modifying files for consistency in R2-4

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# libraryPreparer column (B) for all data rows: fix casing from "J.Plaggenberg" to "J.PLAGGENBERG"
$ws.Range("B2:B42").Value = "J.PLAGGENBERG"

# Update volumePooled (column J) values for consistency across R2-4
$ws.Range("J3").Value = 5
$ws.Range("J6").Value = 5
$ws.Range("J17").Value = 5
$ws.Range("J19").Value = 5
$ws.Range("J20").Value = 5
$ws.Range("J24").Value = 5
$ws.Range("J25").Value = 5
$ws.Range("J27").Value = 0.775
$ws.Range("J28").Value = 0.78
$ws.Range("J33").Value = 5
$ws.Range("J34").Value = 5
$ws.Range("J35").Value = 5
$ws.Range("J36").Value = 5

# Update the active selection to reflect where the author left off editing
$ws.Range("J45").Select()
